$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6130.1
$ws.Range("I70").Value = 22000.5
$ws.Range("J70").Value = 2162.5
$ws.Range("K70").Value = 66001.5
$ws.Range("L70").Value = 6487.5
$ws.Range("M70").Value = -65731.5
$ws.Range("N70").Value = -7027.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6130.1
$ws.Range("I73").Value = 22000.5
$ws.Range("J73").Value = 2162.5
$ws.Range("K73").Value = 66001.5
$ws.Range("L73").Value = 6487.5
$ws.Range("M73").Value = -65065.5
$ws.Range("N73").Value = -8359.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1742.0968
$ws.Range("I113").Value = 1419.9
$ws.Range("J113").Value = 1895.5238
$ws.Range("K113").Value = 1419.9
$ws.Range("L113").Value = 1895.5238
$ws.Range("M113").Value = 1834.1
$ws.Range("N113").Value = -8403.523799999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4959.1665
$ws.Range("I132").Value = 4457.875
$ws.Range("J132").Value = 6964.3335
$ws.Range("K132").Value = 13373.625
$ws.Range("L132").Value = 20893.0005
$ws.Range("M132").Value = -10843.625
$ws.Range("N132").Value = -25953.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 36304.285
$ws.Range("J136").Value = 36304.285
$ws.Range("L136").Value = 36304.285
$ws.Range("N136").Value = -46504.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10642449
$ws.Range("I32").Value = 12661213
$ws.Range("J32").Value = 10287.8
$ws.Range("K32").Value = 12661213
$ws.Range("L32").Value = 10287.8
$ws.Range("M32").Value = -12660926
$ws.Range("N32").Value = -10861.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3629.25
$ws.Range("I74").Value = 3675.4
$ws.Range("J74").Value = 2014
$ws.Range("K74").Value = 3675.4
$ws.Range("L74").Value = 2014
$ws.Range("M74").Value = -2801.4
$ws.Range("N74").Value = -3762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3629.25
$ws.Range("I77").Value = 3675.4
$ws.Range("J77").Value = 2014
$ws.Range("K77").Value = 18377
$ws.Range("L77").Value = 10070
$ws.Range("M77").Value = -14009
$ws.Range("N77").Value = -18806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3095.0667
$ws.Range("I110").Value = 2482.6
$ws.Range("J110").Value = 4320
$ws.Range("K110").Value = 2482.6
$ws.Range("L110").Value = 4320
$ws.Range("M110").Value = -437.5999999999999
$ws.Range("N110").Value = -8410

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2815.0527
$ws.Range("I107").Value = 3037.3333
$ws.Range("J107").Value = 1981.5
$ws.Range("K107").Value = 3037.3333
$ws.Range("L107").Value = 1981.5
$ws.Range("M107").Value = -1117.3333
$ws.Range("N107").Value = -5821.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 816.1818
$ws.Range("I5").Value = 754.2857
$ws.Range("J5").Value = 924.5
$ws.Range("K5").Value = 2262.8571
$ws.Range("L5").Value = 2773.5
$ws.Range("M5").Value = -2150.8571
$ws.Range("N5").Value = -2997.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5318.087
$ws.Range("I131").Value = 608.1667
$ws.Range("J131").Value = 6980.4116
$ws.Range("K131").Value = 1824.5001
$ws.Range("L131").Value = 20941.2348
$ws.Range("M131").Value = 3215.4999
$ws.Range("N131").Value = -31021.2348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 816.1818
$ws.Range("I135").Value = 754.2857
$ws.Range("J135").Value = 924.5
$ws.Range("K135").Value = 6788.571300000001
$ws.Range("L135").Value = 8320.5
$ws.Range("M135").Value = -4253.571300000001
$ws.Range("N135").Value = -13390.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2872.88
$ws.Range("I80").Value = 2726.5334
$ws.Range("J80").Value = 3092.4
$ws.Range("K80").Value = 2726.5334
$ws.Range("L80").Value = 3092.4
$ws.Range("M80").Value = -1728.5334
$ws.Range("N80").Value = -5088.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2872.88
$ws.Range("I83").Value = 2726.5334
$ws.Range("J83").Value = 3092.4
$ws.Range("K83").Value = 13632.667
$ws.Range("L83").Value = 15462
$ws.Range("M83").Value = -8640.666999999999
$ws.Range("N83").Value = -25446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1618
$ws.Range("I97").Value = 990
$ws.Range("J97").Value = 1775
$ws.Range("K97").Value = 990
$ws.Range("L97").Value = 1775
$ws.Range("M97").Value = -494
$ws.Range("N97").Value = -2767

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 14613.786
$ws.Range("I113").Value = 993.6667
$ws.Range("J113").Value = 18328.363
$ws.Range("K113").Value = 993.6667
$ws.Range("L113").Value = 18328.363
$ws.Range("M113").Value = 1176.3333
$ws.Range("N113").Value = -22668.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5909.0835
$ws.Range("I40").Value = 5334
$ws.Range("J40").Value = 6484.1665
$ws.Range("K40").Value = 5334
$ws.Range("L40").Value = 6484.1665
$ws.Range("M40").Value = -5198
$ws.Range("N40").Value = -6756.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1326.9445
$ws.Range("I61").Value = 1034.6364
$ws.Range("J61").Value = 1786.2858
$ws.Range("K61").Value = 1034.6364
$ws.Range("L61").Value = 1786.2858
$ws.Range("M61").Value = -832.6364000000001
$ws.Range("N61").Value = -2190.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2342.75
$ws.Range("I68").Value = 2000.125
$ws.Range("J68").Value = 2685.375
$ws.Range("K68").Value = 2000.125
$ws.Range("L68").Value = 2685.375
$ws.Range("M68").Value = -1251.125
$ws.Range("N68").Value = -4183.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2342.75
$ws.Range("I71").Value = 2000.125
$ws.Range("J71").Value = 2685.375
$ws.Range("K71").Value = 10000.625
$ws.Range("L71").Value = 13426.875
$ws.Range("M71").Value = -6256.625
$ws.Range("N71").Value = -20914.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2518.0588
$ws.Range("I82").Value = 1917.8889
$ws.Range("J82").Value = 3193.25
$ws.Range("K82").Value = 1917.8889
$ws.Range("L82").Value = 3193.25
$ws.Range("M82").Value = -1556.8889
$ws.Range("N82").Value = -3915.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2518.0588
$ws.Range("I85").Value = 1917.8889
$ws.Range("J85").Value = 3193.25
$ws.Range("K85").Value = 1917.8889
$ws.Range("L85").Value = 3193.25
$ws.Range("M85").Value = -669.8888999999999
$ws.Range("N85").Value = -5689.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1326.9445
$ws.Range("I113").Value = 1034.6364
$ws.Range("J113").Value = 1786.2858
$ws.Range("K113").Value = 1034.6364
$ws.Range("L113").Value = 1786.2858
$ws.Range("M113").Value = 1135.3636
$ws.Range("N113").Value = -6126.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4401.1
$ws.Range("I122").Value = 4638
$ws.Range("J122").Value = 4164.2
$ws.Range("K122").Value = 13914
$ws.Range("L122").Value = 12492.6
$ws.Range("M122").Value = -11464
$ws.Range("N122").Value = -17392.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7703954
$ws.Range("I62").Value = 20022380
$ws.Range("J62").Value = 4937.75
$ws.Range("K62").Value = 20022380
$ws.Range("L62").Value = 4937.75
$ws.Range("M62").Value = -20021756
$ws.Range("N62").Value = -6185.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7703954
$ws.Range("I65").Value = 20022380
$ws.Range("J65").Value = 4937.75
$ws.Range("K65").Value = 100111900
$ws.Range("L65").Value = 24688.75
$ws.Range("M65").Value = -100108780
$ws.Range("N65").Value = -30928.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1273.8384
$ws.Range("I132").Value = 1089.9572
$ws.Range("J132").Value = 1717.6897
$ws.Range("K132").Value = 3269.8716
$ws.Range("L132").Value = 5153.0691
$ws.Range("M132").Value = -739.8716000000004
$ws.Range("N132").Value = -10213.0691
